# Revert "new changes in ops (ordercreation & orderpage & order form)"
# Restores the sheet to its prior (pre-"Typist"/"Typist QC" columns, no extra
# "Typing" row) layout: 13 columns (A:M), 3 rows (header + 2 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the extra "Typing" row (row 4) and the stray columns N:O that held
#    data beyond the restored 13-column layout.
# ---------------------------------------------------------------------------
$ws.Rows("4:4").Delete()
$ws.Range("N1:O3").Clear()

# ---------------------------------------------------------------------------
# 2. Rewrite the header row (A1:M1) with the original column order/labels.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Order Received Data and Time"
$ws.Range("B1").Value = "OrderID"
$ws.Range("C1").Value = "Emp ID-Order Assigned"
$ws.Range("D1").Value = "Assignee_QA"
$ws.Range("E1").Value = "Client"
$ws.Range("F1").Value = "Product Name"
$ws.Range("G1").Value = "Lob"
$ws.Range("H1").Value = "Process"
$ws.Range("I1").Value = "State"
$ws.Range("J1").Value = "County"
$ws.Range("K1").Value = "Municipality"
$ws.Range("L1").Value = "Status"
$ws.Range("M1").Value = "Tier"

# ---------------------------------------------------------------------------
# 3. Rewrite data row 2.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 45436
$ws.Range("B2").Value = 1213286
$ws.Range("C2").Value = "SIPL0005"
$ws.Range("D2").Value = "SIPL0004"
$ws.Range("E2").Value = "Reliable Property Reports, Inc"
$ws.Range("F2").Value = "Document Retrieval"
$ws.Range("G2").Value = "Title"
$ws.Range("H2").Value = "Search"
$ws.Range("I2").Value = "FL"
$ws.Range("J2").Value = "Clay"
$ws.Range("K2").Value = "FLClay"
$ws.Range("L2").Value = "WIP"
$ws.Range("M2").Value = "Search(T1) "

# ---------------------------------------------------------------------------
# 4. Rewrite data row 3.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 45439
$ws.Range("B3").Value = 2193289
$ws.Range("C3").Value = "SIPL0005"
$ws.Range("D3").Value = "SIPL0004"
$ws.Range("E3").Value = "Reliable Property Reports, Inc"
$ws.Range("F3").Value = "Marketable Title"
$ws.Range("G3").Value = "Title"
$ws.Range("H3").Value = "Search"
$ws.Range("I3").Value = "FL"
$ws.Range("J3").Value = "Clay"
$ws.Range("K3").Value = "FLClay"
$ws.Range("L3").Value = "WIP"
$ws.Range("M3").Value = "Search(T2)"

# ---------------------------------------------------------------------------
# 5. Formatting: thin black borders everywhere, bold header w/ gold fill,
#    date/time number format on the received-date column.
# ---------------------------------------------------------------------------
$dataRange = $ws.Range("A1:M3")
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

$headerRange = $ws.Range("A1:M1")
$headerRange.Font.Bold = $true
$headerRange.Font.ThemeColor = 1
$headerRange.Interior.Color = 6736351

$ws.Range("A2:A3").NumberFormat = "[$-409]m/d/yy\ h:mm\ AM/PM;@"

# ---------------------------------------------------------------------------
# 6. Column widths matching the restored layout.
# ---------------------------------------------------------------------------
$ws.Columns("C:C").ColumnWidth = 35.4166666666667
$ws.Columns("E:E").ColumnWidth = 11.75
$ws.Columns("F:H").ColumnWidth = 15.0833333333333

# ---------------------------------------------------------------------------
# 7. Selection cursor, as recorded in the reverted workbook.
# ---------------------------------------------------------------------------
$ws.Range("G13").Select()
